$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17 (hunk 0)
$ws.Range("H17").Value = 777.8
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 777.8
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2333.4
$ws.Range("M17").ClearContents() | Out-Null
$ws.Range("N17").Value = -2669.4

# Row 28 (hunk 1)
$ws.Range("H28").Value = 429.1111
$ws.Range("I28").Value = 357.75
$ws.Range("J28").Value = 1000
$ws.Range("K28").Value = 357.75
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = 127.25
$ws.Range("N28").Value = -1970

# Row 75 (hunk 2)
$ws.Range("H75").Value = 67062.8
$ws.Range("J75").Value = 75078.5
$ws.Range("L75").Value = 75078.5
$ws.Range("N75").Value = -76950.5

# Row 78 (hunk 3)
$ws.Range("H78").Value = 67062.8
$ws.Range("J78").Value = 75078.5
$ws.Range("L78").Value = 225235.5
$ws.Range("N78").Value = -234595.5

$ws = $wb.Worksheets.Item("ARM")
# Row 45 (hunk 4)
$ws.Range("H45").Value = 3119.6
$ws.Range("I45").Value = 1599
$ws.Range("K45").Value = 1599
$ws.Range("M45").Value = -1222

# Row 96 (hunk 5)
$ws.Range("H96").Value = 5166.6665
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents() | Out-Null

# Row 101 (hunk 6)
$ws.Range("H101").Value = 19000
$ws.Range("J101").Value = 19000
$ws.Range("L101").Value = 19000
$ws.Range("N101").Value = -25490

# Row 114 (hunk 7)
$ws.Range("H114").Value = 70000
$ws.Range("J114").Value = 70000
$ws.Range("L114").Value = 70000
$ws.Range("N114").Value = -78678

# Row 122 (hunk 8)
$ws.Range("H122").Value = 21968.578
$ws.Range("I122").Value = 18650.25
$ws.Range("K122").Value = 55950.75
$ws.Range("M122").Value = -53500.75

$ws = $wb.Worksheets.Item("BSM")
# Row 106 (hunk 9)
$ws.Range("H106").Value = 28500
$ws.Range("J106").Value = 28500
$ws.Range("L106").Value = 28500
$ws.Range("N106").Value = -31024

# Row 128 (hunk 10)
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("CRP")
# Row 58 (hunk 11)
$ws.Range("H58").Value = 3909.5
$ws.Range("I58").Value = 1970.6666
$ws.Range("J58").Value = 5848.3335
$ws.Range("K58").Value = 1970.6666
$ws.Range("L58").Value = 5848.3335
$ws.Range("M58").Value = -1767.6666
$ws.Range("N58").Value = -6254.3335

# Row 62 (hunk 12)
$ws.Range("H62").Value = 399999
$ws.Range("J62").Value = 399999
$ws.Range("L62").Value = 399999
$ws.Range("N62").Value = -401247

# Row 65 (hunk 13)
$ws.Range("H65").Value = 399999
$ws.Range("J65").Value = 399999
$ws.Range("L65").Value = 1999995
$ws.Range("N65").Value = -2006235

# Row 86 (hunk 14)
$ws.Range("H86").Value = 3274.75

# Row 89 (hunk 15)
$ws.Range("H89").Value = 3274.75

# Row 136 (hunk 16)
$ws.Range("H136").Value = 3909.5
$ws.Range("I136").Value = 1970.6666
$ws.Range("J136").Value = 5848.3335
$ws.Range("K136").Value = 5911.9998
$ws.Range("L136").Value = 17545.0005
$ws.Range("M136").Value = -3361.9998
$ws.Range("N136").Value = -22645.0005

$ws = $wb.Worksheets.Item("CUL")
# Row 7 (hunk 17)
$ws.Range("H7").Value = 127.333336
$ws.Range("J7").Value = 150
$ws.Range("L7").Value = 450
$ws.Range("N7").Value = -674

# Row 21 (hunk 18)
$ws.Range("H21").Value = 140.66667
$ws.Range("J21").Value = 125
$ws.Range("L21").Value = 375
$ws.Range("N21").Value = -721

# Row 98 (hunk 19)
$ws.Range("H98").Value = 2232.2222
$ws.Range("I98").Value = 2749.5
$ws.Range("J98").Value = 2084.4285
$ws.Range("K98").Value = 8248.5
$ws.Range("L98").Value = 6253.2855
$ws.Range("M98").Value = -6750.5
$ws.Range("N98").Value = -9249.2855

# Row 107 (hunk 20)
$ws.Range("H107").Value = 143574.72
$ws.Range("J107").Value = 143574.72
$ws.Range("L107").Value = 430724.16
$ws.Range("N107").Value = -434564.16

# Row 121 (hunk 21)
$ws.Range("H121").Value = 765.4286
$ws.Range("I121").Value = 365
$ws.Range("J121").Value = 1766.5
$ws.Range("K121").Value = 1095
$ws.Range("L121").Value = 5299.5
$ws.Range("M121").Value = 215
$ws.Range("N121").Value = -7919.5

# Row 129 (hunk 22)
$ws.Range("H129").Value = 1182.5
$ws.Range("I129").Value = 492
$ws.Range("J129").Value = 2333.3333
$ws.Range("K129").Value = 1476
$ws.Range("L129").Value = 6999.999899999999
$ws.Range("M129").Value = 3524
$ws.Range("N129").Value = -16999.9999

# Row 131 (hunk 23)
$ws.Range("H131").Value = 1381.1063
$ws.Range("J131").Value = 1404.0698
$ws.Range("L131").Value = 4212.2094
$ws.Range("N131").Value = -14292.2094

$ws = $wb.Worksheets.Item("GSM")
# Row 59 (hunk 24)
$ws.Range("H59").Value = 1000
$ws.Range("J59").Value = 1000
$ws.Range("L59").Value = 1000
$ws.Range("N59").Value = -2166

# Row 80 (hunk 25)
$ws.Range("H80").Value = 3975.8
$ws.Range("I80").Value = 2959.6667
$ws.Range("K80").Value = 2959.6667
$ws.Range("M80").Value = -1961.6667

# Row 83 (hunk 26)
$ws.Range("H83").Value = 3975.8
$ws.Range("I83").Value = 2959.6667
$ws.Range("K83").Value = 14798.3335
$ws.Range("M83").Value = -9806.333500000001

# Row 122 (hunk 27)
$ws.Range("H122").Value = 113817.78
$ws.Range("I122").Value = 1872.4
$ws.Range("J122").Value = 253749.5
$ws.Range("K122").Value = 5617.200000000001
$ws.Range("L122").Value = 761248.5
$ws.Range("M122").Value = -3167.200000000001
$ws.Range("N122").Value = -766148.5

# Row 131 (hunk 28)
$ws.Range("H131").Value = 20000
$ws.Range("J131").Value = 20000
$ws.Range("L131").Value = 20000
$ws.Range("N131").Value = -30080

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (hunk 29)
$ws.Range("H7").Value = 4615.3335
$ws.Range("I7").Value = 3931.6667
$ws.Range("K7").Value = 3931.6667
$ws.Range("M7").Value = -3819.6667

# Row 40 (hunk 30)
$ws.Range("H40").Value = 5039.25
$ws.Range("I40").Value = 4897.1
$ws.Range("K40").Value = 4897.1
$ws.Range("M40").Value = -4761.1

# Row 122 (hunk 31)
$ws.Range("H122").Value = 4002
$ws.Range("I122").Value = 4002
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 12006
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -9556
$ws.Range("N122").ClearContents() | Out-Null

# Row 124 (hunk 32)
$ws.Range("H124").Value = 34606.75
$ws.Range("J124").Value = 34606.75
$ws.Range("L124").Value = 34606.75
$ws.Range("N124").Value = -44426.75

# Row 126 (hunk 33)
$ws.Range("H126").Value = 4615.3335
$ws.Range("I126").Value = 3931.6667
$ws.Range("K126").Value = 11795.0001
$ws.Range("M126").Value = -9325.000100000001

$ws = $wb.Worksheets.Item("WVR")
# Row 3 (hunk 34)
$ws.Range("H3").Value = 30222.555
$ws.Range("I3").Value = 52921.2
$ws.Range("J3").Value = 1849.25
$ws.Range("K3").Value = 52921.2
$ws.Range("L3").Value = 1849.25
$ws.Range("M3").Value = -52807.2
$ws.Range("N3").Value = -2077.25

# Row 101 (hunk 35)
$ws.Range("H101").Value = 1600
$ws.Range("J101").Value = 1600
$ws.Range("L101").Value = 1600
$ws.Range("N101").Value = -8090

# Row 103 (hunk 36)
$ws.Range("H103").Value = 19999.5
$ws.Range("J103").Value = 19999.5
$ws.Range("L103").Value = 19999.5
$ws.Range("N103").Value = -22343.5

# Row 107 (hunk 37)
$ws.Range("H107").Value = 500.66666
$ws.Range("I107").Value = 500.66666
$ws.Range("K107").Value = 1501.99998
$ws.Range("M107").Value = 418.0000199999999
